$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A1 previously held "foo"; it now holds "Sample Name".
# B1 previously held "bar"; it still holds "bar" (unchanged content).
$ws.Range("A1").Value = "Sample Name"
$ws.Range("B1").Value = "bar"

# Selection moves to B4.
$ws.Range("B4").Select()
